$wb = $excel.ActiveWorkbook

$zongji = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# --- 1. Create the new "2022-Q4" sheet ---
# Duplicate the "2022-Q3" sheet (same columns / fund identifiers in A-C) and
# place the copy right after "总计", i.e. before "2022-Q3". This gives us a
# correctly-typed template (fund codes/names) that we then overwrite with the
# new quarter's figures.
$q3.Copy($null, $zongji)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$q4.Range("D2").Value = "'14.15"
$q4.Range("E2").Value = "'75.21"
$q4.Range("F2").Value = "'4.08"
$q4.Range("G2").Value = "'0.5773"
$q4.Range("H2").Value = 4

$q4.Range("D3").Value = "'0.39"
$q4.Range("E3").Value = "'75.21"
$q4.Range("F3").Value = "'4.08"
$q4.Range("G3").Value = "'0.0159"
$q4.Range("H3").Value = 4

# --- 2. Add the 2022-Q4 row to the "总计" summary sheet ---
# Copy the style of the last existing row down into a new row 6 so the new
# row's formatting matches the others, then fill in the 2021-Q1 data that
# now belongs there (everything shifts down by one row).
$zongji.Range("A5").Copy($zongji.Range("A6"))
$zongji.Range("B5").Copy($zongji.Range("B6"))
$zongji.Range("C5").Copy($zongji.Range("C6"))
$zongji.Range("D5").Copy($zongji.Range("D6"))

$zongji.Range("A6").Value = 4
$zongji.Range("B6").Value = "2021-Q1"
$zongji.Range("C6").Value = 2
$zongji.Range("D6").Value = 0.01

$zongji.Range("A5").Value = 3
$zongji.Range("B5").Value = "2022-Q1"
$zongji.Range("C5").Value = 2
$zongji.Range("D5").Value = 0.78

$zongji.Range("A4").Value = 2
$zongji.Range("B4").Value = "2022-Q2"
$zongji.Range("C4").Value = 2
$zongji.Range("D4").Value = 0.71

$zongji.Range("A3").Value = 1
$zongji.Range("B3").Value = "2022-Q3"
$zongji.Range("C3").Value = 2
$zongji.Range("D3").Value = 0.7

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q4"
$zongji.Range("C2").Value = 2
$zongji.Range("D2").Value = 0.59

# Restore the originally selected tab ("2021-Q1" was the active sheet).
$wb.Worksheets.Item("2021-Q1").Activate()
